# Applies the cryptos.xlsx price/volume/coin update described in the commit
# "Updated cryptos list on Thu Feb 29 09:20:45 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '62.722.36'
$ws.Range("E2").Value = '  +5.92%  '

# Row 3
$ws.Range("D3").Value = '3.474.91'
$ws.Range("E3").Value = '  +4.50%  '

# Row 4
$ws.Range("E4").Value = '  +0.26%  '

# Row 5
$ws.Range("D5").Value = "'417.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.57%  '

# Row 6
$ws.Range("D6").Value = "'131.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +18.79%  '

# Row 7
$ws.Range("D7").Value = '3.462.52'
$ws.Range("E7").Value = '  +4.45%  '

# Row 8
$ws.Range("D8").Value = "'0.595"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.61%  '

# Row 9
$ws.Range("E9").Value = '  +0.13%  '

# Row 10
$ws.Range("D10").Value = "'0.692"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.85%  '

# Row 11
$ws.Range("E11").Value = '  +28.58%  '

# Row 12
$ws.Range("D12").Value = "'44.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +11.17%  '

# Row 13
$ws.Range("E13").Value = '  +0.48%  '

# Row 14
$ws.Range("D14").Value = '4.034.73'
$ws.Range("E14").Value = '  +5.03%  '

# Row 15
$ws.Range("D15").Value = "'8.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.74%  '

# Row 16
$ws.Range("D16").Value = "'20.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.68%  '

# Row 17
$ws.Range("D17").Value = '3.481.27'
$ws.Range("E17").Value = '  +4.27%  '

# Row 18
$ws.Range("B18").Value = 'Polygon'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D18").Value = "'1.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.66%  '

# Row 19
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '62.803.18'
$ws.Range("E19").Value = '  +6.21%  '

# Row 20
$ws.Range("E20").Value = '  +3.20%  '

# Row 21
$ws.Range("D21").Value = "'0.0000136"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +25.80%  '

# Row 22
$ws.Range("D22").Value = "'3.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.36%  '

# Row 23
$ws.Range("D23").Value = "'13.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.99%  '

# Row 24
$ws.Range("D24").Value = "'82.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +9.32%  '

# Row 25
$ws.Range("D25").Value = "'315.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.15%  '

# Row 26
$ws.Range("D26").Value = "'3.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.18%  '

# Row 27
$ws.Range("D27").Value = "'30.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.40%  '

# Row 28
$ws.Range("D28").Value = "'8.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.21%  '

# Row 29
$ws.Range("D29").Value = "'7.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.19%  '

# Row 30
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = "'0.123"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +9.91%  '

# Row 31
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").Value = "'0.179"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.79%  '

# Row 32
$ws.Range("D32").Value = "'4.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.43%  '

# Row 33
$ws.Range("D33").Value = "'44.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +10.94%  '

# Row 34
$ws.Range("D34").Value = "'11.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.37%  '

# Row 35
$ws.Range("D35").Value = "'2.64"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +24.04%  '

# Row 36
$ws.Range("E36").Value = '  +0.19%  '

# Row 37
$ws.Range("E37").Value = '  -4.87%  '

# Row 38
$ws.Range("D38").Value = "'52.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.47%  '

# Row 39
$ws.Range("D39").Value = "'3.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.93%  '

# Row 40
$ws.Range("D40").Value = "'0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.10%  '

# Row 41
$ws.Range("D41").Value = "'3.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -6.77%  '

# Row 42
$ws.Range("D42").Value = "'2.02"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.36%  '

# Row 43
$ws.Range("E43").Value = '  +3.25%  '

# Row 44
$ws.Range("D44").Value = "'137.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.21%  '

# Row 45
$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").Value = "'17.61"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.70%  '

# Row 46
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").Value = "'4.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.71%  '

# Row 47
$ws.Range("D47").Value = "'0.289"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.61%  '

# Row 48
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = "'2.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.51%  '

# Row 49
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = "'22.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.94%  '

# Row 50
$ws.Range("D50").Value = '2.256.36'

# Row 51
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D51").Value = "'2.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.62%  '
